$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It belongs right
# after the existing row 65 (chronologically it is just another entry in
# the same ongoing series), so insert a fresh row at position 66 and
# shift every following record (old rows 66-97) down by one, ending up
# at rows 67-98.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new record's data. All of
# the "constant" columns (market, region, category codes/names, unit of
# sale, origin, classification, etc.) repeat the same values used
# throughout this sheet.
$ws.Cells.Item(66, 1).Value = 3
$ws.Cells.Item(66, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(66, 3).Value = "Coquimbo"
$ws.Cells.Item(66, 4).Value = 44845
$ws.Cells.Item(66, 5).Value = 5
$ws.Cells.Item(66, 6).Value = 100112035
$ws.Cells.Item(66, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 42
$ws.Cells.Item(66, 11).Value = 15000
$ws.Cells.Item(66, 12).Value = 15000
$ws.Cells.Item(66, 13).Value = 15000
$ws.Cells.Item(66, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(66, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(66, 16).Value = 1000
$ws.Cells.Item(66, 17).Value = 15
$ws.Cells.Item(66, 18).Value = "Hortaliza"
